$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 478.09525
$ws.Range("I28").Value = 370.73685
$ws.Range("K28").Value = 370.73685
$ws.Range("M28").Value = 114.26315
$ws.Range("H55").Value = 144.41176
$ws.Range("I55").Value = 144.16667
$ws.Range("J55").Value = 145
$ws.Range("K55").Value = 144.16667
$ws.Range("L55").Value = 145
$ws.Range("M55").Value = 69.83332999999999
$ws.Range("N55").Value = -573
$ws.Range("H62").Value = 16256
$ws.Range("I62").Value = 12388.048
$ws.Range("K62").Value = 12388.048
$ws.Range("M62").Value = -11764.048
$ws.Range("H65").Value = 16256
$ws.Range("I65").Value = 12388.048
$ws.Range("K65").Value = 61940.24000000001
$ws.Range("M65").Value = -58820.24000000001
$ws.Range("H86").Value = 2664.4443
$ws.Range("I86").Value = 1678.75
$ws.Range("J86").Value = 4635.8335
$ws.Range("K86").Value = 1678.75
$ws.Range("L86").Value = 4635.8335
$ws.Range("M86").Value = -555.75
$ws.Range("N86").Value = -6881.8335
$ws.Range("H89").Value = 2664.4443
$ws.Range("I89").Value = 1678.75
$ws.Range("J89").Value = 4635.8335
$ws.Range("K89").Value = 8393.75
$ws.Range("L89").Value = 23179.1675
$ws.Range("M89").Value = -2777.75
$ws.Range("N89").Value = -34411.1675
$ws.Range("H98").Value = 2218.6428
$ws.Range("I98").Value = 2218.6428
$ws.Range("K98").Value = 2218.6428
$ws.Range("M98").Value = -720.6428000000001
$ws.Range("H112").Value = 1410.0294
$ws.Range("J112").Value = 1529.5186
$ws.Range("L112").Value = 4588.5558
$ws.Range("N112").Value = -6804.5558
$ws.Range("H122").Value = 2218.6428
$ws.Range("I122").Value = 2218.6428
$ws.Range("K122").Value = 6655.928400000001
$ws.Range("M122").Value = -4205.928400000001
$ws.Range("H137").Value = 2690.5173
$ws.Range("I137").Value = 1827
$ws.Range("J137").Value = 2915.7827
$ws.Range("K137").Value = 5481
$ws.Range("L137").Value = 8747.348100000001
$ws.Range("M137").Value = -2931
$ws.Range("N137").Value = -13847.3481
$ws.Range("H138").Value = 2744
$ws.Range("I138").Value = 1780.9166
$ws.Range("J138").Value = 3748.9565
$ws.Range("K138").Value = 5342.7498
$ws.Range("L138").Value = 11246.8695
$ws.Range("M138").Value = -202.7497999999996
$ws.Range("N138").Value = -21526.8695
$ws.Range("H141").Value = 903.875
$ws.Range("I141").Value = 936.2857
$ws.Range("J141").Value = 677
$ws.Range("K141").Value = 2808.8571
$ws.Range("L141").Value = 2031
$ws.Range("M141").Value = 2371.1429
$ws.Range("N141").Value = -12391

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3628.8333
$ws.Range("I45").Value = 3628.8333
$ws.Range("K45").Value = 3628.8333
$ws.Range("M45").Value = -3251.8333
$ws.Range("H61").Value = 8811.894
$ws.Range("I61").Value = 8683.625
$ws.Range("J61").Value = 8945.739
$ws.Range("K61").Value = 8683.625
$ws.Range("L61").Value = 8945.739
$ws.Range("M61").Value = -8471.625
$ws.Range("N61").Value = -9369.739
$ws.Range("H74").Value = 3460.0635
$ws.Range("I74").Value = 2971.1702
$ws.Range("J74").Value = 4896.1875
$ws.Range("K74").Value = 2971.1702
$ws.Range("L74").Value = 4896.1875
$ws.Range("M74").Value = -2097.1702
$ws.Range("N74").Value = -6644.1875
$ws.Range("H77").Value = 3460.0635
$ws.Range("I77").Value = 2971.1702
$ws.Range("J77").Value = 4896.1875
$ws.Range("K77").Value = 14855.851
$ws.Range("L77").Value = 24480.9375
$ws.Range("M77").Value = -10487.851
$ws.Range("N77").Value = -33216.9375
$ws.Range("H132").Value = 6753.271
$ws.Range("I132").Value = 5861.3076
$ws.Range("K132").Value = 17583.9228
$ws.Range("M132").Value = -15053.9228
$ws.Range("H135").Value = 44949.5
$ws.Range("I135").Value = 45000
$ws.Range("J135").Value = 44899
$ws.Range("K135").Value = 45000
$ws.Range("L135").Value = 44899
$ws.Range("M135").Value = -39930
$ws.Range("N135").Value = -55039
$ws.Range("H136").Value = 8811.894
$ws.Range("I136").Value = 8683.625
$ws.Range("J136").Value = 8945.739
$ws.Range("K136").Value = 26050.875
$ws.Range("L136").Value = 26837.217
$ws.Range("M136").Value = -23500.875
$ws.Range("N136").Value = -31937.217
$ws.Range("H139").Value = 64500
$ws.Range("J139").Value = 64500
$ws.Range("L139").Value = 64500
$ws.Range("N139").Value = -74780

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1044.75
$ws.Range("I36").Value = 1044.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1044.75
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -510.75
$ws.Range("H86").Value = 1290.2858
$ws.Range("I86").Value = 1275.4615
$ws.Range("K86").Value = 1275.4615
$ws.Range("M86").Value = -152.4614999999999
$ws.Range("H89").Value = 1290.2858
$ws.Range("I89").Value = 1275.4615
$ws.Range("K89").Value = 6377.307499999999
$ws.Range("M89").Value = -761.307499999999
$ws.Range("H94").Value = 710.2308
$ws.Range("I94").Value = 710.2308
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 710.2308
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -259.2308
$ws.Range("H134").Value = 16933
$ws.Range("I134").Value = 20812.51
$ws.Range("K134").Value = 62437.53
$ws.Range("M134").Value = -59902.53

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6827.549
$ws.Range("I31").Value = 1675.6818
$ws.Range("J31").Value = 10735.862
$ws.Range("K31").Value = 1675.6818
$ws.Range("L31").Value = 10735.862
$ws.Range("M31").Value = -1380.6818
$ws.Range("N31").Value = -11325.862
$ws.Range("H34").Value = 6827.549
$ws.Range("I34").Value = 1675.6818
$ws.Range("J34").Value = 10735.862
$ws.Range("K34").Value = 1675.6818
$ws.Range("L34").Value = 10735.862
$ws.Range("M34").Value = -1473.6818
$ws.Range("N34").Value = -11139.862
$ws.Range("H134").Value = 22604
$ws.Range("I134").Value = 27255.643
$ws.Range("K134").Value = 81766.929
$ws.Range("M134").Value = -79231.929

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 17505.5
$ws.Range("H90").Value = 17505.5
$ws.Range("H118").Value = 1464.4
$ws.Range("I118").Value = 1480.5
$ws.Range("K118").Value = 4441.5
$ws.Range("M118").Value = -3198.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.65218
$ws.Range("J2").Value = 323.5
$ws.Range("L2").Value = 323.5
$ws.Range("N2").Value = -549.5
$ws.Range("H122").Value = 4623.3477
$ws.Range("I122").Value = 2309.5625
$ws.Range("K122").Value = 6928.6875
$ws.Range("M122").Value = -4478.6875
$ws.Range("H132").Value = 2489.182
$ws.Range("I132").Value = 2433.4119
$ws.Range("J132").Value = 2678.8
$ws.Range("K132").Value = 7300.2357
$ws.Range("L132").Value = 8036.400000000001
$ws.Range("M132").Value = -4770.2357
$ws.Range("N132").Value = -13096.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5946.357
$ws.Range("J68").Value = 6619.048
$ws.Range("L68").Value = 6619.048
$ws.Range("N68").Value = -8117.048
$ws.Range("H71").Value = 5946.357
$ws.Range("J71").Value = 6619.048
$ws.Range("L71").Value = 33095.24
$ws.Range("N71").Value = -40583.24
$ws.Range("H93").Value = 2648.6667
$ws.Range("I93").Value = 2780.3635
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 2780.3635
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = -1532.3635
$ws.Range("N93").Value = -3696
$ws.Range("H122").Value = 3501.48
$ws.Range("I122").Value = 3236.1904
$ws.Range("K122").Value = 9708.5712
$ws.Range("M122").Value = -7258.5712
$ws.Range("H132").Value = 5070.3145
$ws.Range("I132").Value = 4281
$ws.Range("K132").Value = 12843
$ws.Range("M132").Value = -10313
$ws.Range("H136").Value = 7287.8184
$ws.Range("I136").Value = 6947.1763
$ws.Range("J136").Value = 8446
$ws.Range("K136").Value = 20841.5289
$ws.Range("L136").Value = 25338
$ws.Range("M136").Value = -18291.5289
$ws.Range("N136").Value = -30438

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2107.8
$ws.Range("I81").Value = 2160
$ws.Range("J81").Value = 1899
$ws.Range("K81").Value = 4320
$ws.Range("L81").Value = 3798
$ws.Range("M81").Value = -3259
$ws.Range("N81").Value = -5920
$ws.Range("H84").Value = 2107.8
$ws.Range("I84").Value = 2160
$ws.Range("J84").Value = 1899
$ws.Range("K84").Value = 21600
$ws.Range("L84").Value = 18990
$ws.Range("M84").Value = -16296
$ws.Range("N84").Value = -29598
$ws.Range("H132").Value = 15726819
$ws.Range("I132").Value = 8052.6
$ws.Range("J132").Value = 22276304
$ws.Range("K132").Value = 24157.8
$ws.Range("L132").Value = 66828912
$ws.Range("M132").Value = -21627.8
$ws.Range("N132").Value = -66833972
$ws.Range("H136").Value = 7374718.5
$ws.Range("I136").Value = 8060308.5
$ws.Range("K136").Value = 24180925.5
$ws.Range("M136").Value = -24178375.5
